# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullets as impact-focused
# accomplishment statements (per commit message), reducing the list from
# six job-duty style bullets to four concise achievement statements.

$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading so we operate on the
# correct section. (Similar/identical bullet text also appears earlier in
# the document, under PROFESSIONAL EXPERIENCE, and must be left untouched.)
$achievementsHeadingRange = $d.Content.Duplicate
$achievementsHeadingRange.Find.ClearFormatting()
$found = $achievementsHeadingRange.Find.Execute("KEY ACHIEVEMENTS AND IMPACT", `
    $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading"
}
$headingStart = $achievementsHeadingRange.Paragraphs(1).Range.Start

# Walk forward from the heading to find the six bullet ("•") paragraphs
# that currently sit underneath the "Impact" sub-heading.
$paraCount = $d.Paragraphs.Count
$bulletIndices = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Start -gt $headingStart -and $para.Range.Text.TrimStart().StartsWith("•")) {
        [void]$bulletIndices.Add($i)
        if ($bulletIndices.Count -eq 6) {
            break
        }
    }
}

if ($bulletIndices.Count -ne 6) {
    throw "Expected 6 bullet paragraphs under KEY ACHIEVEMENTS AND IMPACT, found $($bulletIndices.Count)"
}

# Replace the whole span (first bullet paragraph through the last) in a
# single assignment, using a carriage-return to separate each new bullet.
# This rewrites the six old job-duty bullets into four new, concise,
# impact-focused accomplishment statements.
$firstBullet = $d.Paragraphs($bulletIndices[0])
$lastBullet = $d.Paragraphs($bulletIndices[5])
$bulletsRange = $d.Range($firstBullet.Range.Start, $lastBullet.Range.End)

$newBullets = @(
    "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%",
    "• `$4.7M savings enabled nonprofit access",
    "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions",
    "• 178% accuracy improvement in racial classification algorithms"
) -join "`r"

$bulletsRange.Text = $newBullets
